$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.167.19"
$ws.Range("E2").Value = "  -0.31%  "
$ws.Range("D3").Value = "3.302.10"
$ws.Range("E3").Value = "  -1.63%  "
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "189.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.75%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "559.89"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D9").Value = "3.296.51"
$ws.Range("E9").Value = "  -1.64%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.184"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.68%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.587"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.98%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "47.64"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.47%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000271"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.95%  "
$ws.Range("E14").Value = "  -0.63%  "
$ws.Range("D15").Value = "3.829.89"
$ws.Range("E15").Value = "  -1.73%  "
$ws.Range("E16").Value = "  -0.40%  "
$ws.Range("D17").Value = "66.130.99"
$ws.Range("E17").Value = "  -0.32%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.04"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.83%  "
$ws.Range("E19").Value = "  +0.28%  "
$ws.Range("D20").Value = "3.297.76"
$ws.Range("E20").Value = "  -1.92%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.09"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.40%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.912"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "18.43"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +9.74%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.12"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "100.80"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.73%  "
$ws.Range("E26").Value = "  -2.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.16%  "
$ws.Range("E28").Value = "  +1.31%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.69"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.34%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.65"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.61%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "30.28"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.14%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.72"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.68%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.07"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.18%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "567.85"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.57%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "11.10"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.18%  "
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "57.24"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.44%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.07%  "
$ws.Range("D39").Value = "3.706.52"
$ws.Range("E39").Value = "  -2.94%  "
$ws.Range("D40").Value = "0.0₃0727"
$ws.Range("E40").Value = "  +1.09%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "33.96"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.39%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.30"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.89%  "
$ws.Range("E43").Value = "  +1.98%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.71"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.58%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.40"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.33%  "
$ws.Range("E46").Value = "  -1.86%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0422"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.60%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.27"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.13%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.130"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.50%  "
$ws.Range("E50").Value = "  -2.13%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.999"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.17%  "
